$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 91 (shifts old rows 91-103 down to 93-105)
$ws.Rows("91:92").Insert()

# New row 91 data
$ws.Cells.Item(91,1).Value = 9
$ws.Cells.Item(91,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(91,3).Value = "Metropolitana"
$ws.Cells.Item(91,4).Value = 44474
$ws.Cells.Item(91,5).Value = 13
$ws.Cells.Item(91,6).Value = "Fruta"
$ws.Cells.Item(91,7).Value = 100101
$ws.Cells.Item(91,8).Value = "Berries"
$ws.Cells.Item(91,9).Value = 100101001
$ws.Cells.Item(91,10).Value = "Arándano (blue)"
$ws.Cells.Item(91,11).Value = "Sin especificar"
$ws.Cells.Item(91,12).Value = "Primera"
$ws.Cells.Item(91,13).Value = 380
$ws.Cells.Item(91,14).Value = 12000
$ws.Cells.Item(91,15).Value = 12000
$ws.Cells.Item(91,16).Value = 12000
$ws.Cells.Item(91,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(91,18).Value = "Provincia de Linares"
$ws.Cells.Item(91,19).Value = 6000
$ws.Cells.Item(91,20).Value = 2

# New row 92 data
$ws.Cells.Item(92,1).Value = 9
$ws.Cells.Item(92,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(92,3).Value = "Metropolitana"
$ws.Cells.Item(92,4).Value = 44474
$ws.Cells.Item(92,5).Value = 13
$ws.Cells.Item(92,6).Value = "Fruta"
$ws.Cells.Item(92,7).Value = 100101
$ws.Cells.Item(92,8).Value = "Berries"
$ws.Cells.Item(92,9).Value = 100101001
$ws.Cells.Item(92,10).Value = "Arándano (blue)"
$ws.Cells.Item(92,11).Value = "Sin especificar"
$ws.Cells.Item(92,12).Value = "Segunda"
$ws.Cells.Item(92,13).Value = 450
$ws.Cells.Item(92,14).Value = 10000
$ws.Cells.Item(92,15).Value = 10000
$ws.Cells.Item(92,16).Value = 10000
$ws.Cells.Item(92,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(92,18).Value = "Provincia de Linares"
$ws.Cells.Item(92,19).Value = 5000
$ws.Cells.Item(92,20).Value = 2
